# Auto-generated edit script: updates Leve profit-tracking cells per scheduled market-data refresh.
# Applies per-cell numeric updates (and clears cells that no longer have data) across 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1525.2106
$ws.Range("I2").Value = 1425.6774
$ws.Range("J2").Value = 1966
$ws.Range("K2").Value = 1425.6774
$ws.Range("L2").Value = 1966
$ws.Range("M2").Value = -1312.6774
$ws.Range("N2").Value = -2192
$ws.Range("H74").Value = 45455420
$ws.Range("I74").Value = 76923736
$ws.Range("J74").Value = 1186.4445
$ws.Range("K74").Value = 76923736
$ws.Range("L74").Value = 1186.4445
$ws.Range("M74").Value = -76922862
$ws.Range("N74").Value = -2934.4445
$ws.Range("H77").Value = 45455420
$ws.Range("I77").Value = 76923736
$ws.Range("J77").Value = 1186.4445
$ws.Range("K77").Value = 384618680
$ws.Range("L77").Value = 5932.2225
$ws.Range("M77").Value = -384614312
$ws.Range("N77").Value = -14668.2225
$ws.Range("H102").Value = 1526.25
$ws.Range("I102").Value = 1077.5
$ws.Range("K102").Value = 1077.5
$ws.Range("M102").Value = 544.5
$ws.Range("H116").Value = 1525.2106
$ws.Range("I116").Value = 1425.6774
$ws.Range("J116").Value = 1966
$ws.Range("K116").Value = 1425.6774
$ws.Range("L116").Value = 1966
$ws.Range("M116").Value = 868.3226
$ws.Range("N116").Value = -6554
$ws.Range("H122").Value = 3242.6875
$ws.Range("I122").Value = 2835.0908
$ws.Range("J122").Value = 4139.4
$ws.Range("K122").Value = 8505.2724
$ws.Range("L122").Value = 12418.2
$ws.Range("M122").Value = -6055.2724
$ws.Range("N122").Value = -17318.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1525.2106
$ws.Range("I3").Value = 1425.6774
$ws.Range("J3").Value = 1966
$ws.Range("K3").Value = 1425.6774
$ws.Range("L3").Value = 1966
$ws.Range("M3").Value = -1311.6774
$ws.Range("N3").Value = -2194

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13682.767
$ws.Range("I31").Value = 23948.215
$ws.Range("K31").Value = 23948.215
$ws.Range("M31").Value = -23653.215
$ws.Range("H34").Value = 13682.767
$ws.Range("I34").Value = 23948.215
$ws.Range("K34").Value = 23948.215
$ws.Range("M34").Value = -23746.215
$ws.Range("H62").Value = 7003
$ws.Range("I62").Value = 7500
$ws.Range("J62").Value = 6506
$ws.Range("K62").Value = 7500
$ws.Range("L62").Value = 6506
$ws.Range("M62").Value = -6876
$ws.Range("N62").Value = -7754
$ws.Range("H65").Value = 7003
$ws.Range("I65").Value = 7500
$ws.Range("J65").Value = 6506
$ws.Range("K65").Value = 37500
$ws.Range("L65").Value = 32530
$ws.Range("M65").Value = -34380
$ws.Range("N65").Value = -38770

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 754.5700000000001
$ws.Range("I131").Value = 25
$ws.Range("J131").Value = 769.45917
$ws.Range("K131").Value = 75
$ws.Range("L131").Value = 2308.37751
$ws.Range("M131").Value = 4965
$ws.Range("N131").Value = -12388.37751
$ws.Range("H134").Value = 3594.2632
$ws.Range("I134").Value = 1293
$ws.Range("J134").Value = 6151.222
$ws.Range("K134").Value = 3879
$ws.Range("L134").Value = 18453.666
$ws.Range("M134").Value = 1191
$ws.Range("N134").Value = -28593.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4719.9
$ws.Range("I22").Value = 3814.2856
$ws.Range("J22").Value = 6833
$ws.Range("K22").Value = 3814.2856
$ws.Range("L22").Value = 6833
$ws.Range("M22").Value = -3519.2856
$ws.Range("N22").Value = -7423
$ws.Range("H27").Value = 4719.9
$ws.Range("I27").Value = 3814.2856
$ws.Range("J27").Value = 6833
$ws.Range("K27").Value = 3814.2856
$ws.Range("L27").Value = 6833
$ws.Range("M27").Value = -3707.2856
$ws.Range("N27").Value = -7047

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 35123.668
$ws.Range("J27").Value = 35123.668
$ws.Range("L27").Value = 35123.668
$ws.Range("N27").Value = -35261.668
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16040
$ws.Range("H62").Value = 5687.875
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5687.875
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5687.875
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6935.875
$ws.Range("H65").Value = 5687.875
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5687.875
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 28439.375
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -34679.375
$ws.Range("H81").Value = 1949.8334
$ws.Range("I81").Value = 1739.8
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 3479.6
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -2418.6
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 1949.8334
$ws.Range("I84").Value = 1739.8
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 17398
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -12094
$ws.Range("N84").Value = -40608
$ws.Range("H113").Value = 1042.3684
$ws.Range("I113").Value = 1489.8334
$ws.Range("J113").Value = 275.2857
$ws.Range("K113").Value = 4469.5002
$ws.Range("L113").Value = 825.8571000000001
$ws.Range("M113").Value = -2299.5002
$ws.Range("N113").Value = -5165.8571
$ws.Range("H122").Value = 1780.5385
$ws.Range("I122").Value = 1758.2858
$ws.Range("K122").Value = 5274.857400000001
$ws.Range("M122").Value = -2824.857400000001
$ws.Range("H126").Value = 1200
$ws.Range("I126").Value = 1200
$ws.Range("K126").Value = 3600
$ws.Range("M126").Value = -1130
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 906.4091
$ws.Range("I132").Value = 635.1142599999999
$ws.Range("J132").Value = 1961.4445
$ws.Range("K132").Value = 1905.34278
$ws.Range("L132").Value = 5884.333500000001
$ws.Range("M132").Value = 624.6572200000001
$ws.Range("N132").Value = -10944.3335
$ws.Range("H135").Value = 50615
$ws.Range("J135").Value = 50615
$ws.Range("L135").Value = 50615
$ws.Range("N135").Value = -60755
$ws.Range("H137").Value = 47715
$ws.Range("J137").Value = 47715
$ws.Range("L137").Value = 47715
$ws.Range("N137").Value = -57915
